# Notas_Master.xlsx edit script
# - Sheet "803": update several "Heteroevaluación" grades (column H)
# - Sheet "115": update several "Guía 1 de estadistica" grades (column H),
#   and insert a brand-new 9-row grade block for student
#   OSPINA ESCOBAR, SEBASTIAN (matricula 240667 / doc 1035974995) right
#   before the existing PLANCHEZ URDANETA block (which starts at row 164),
#   pushing every later row down by 9.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "803" — single-cell grade corrections
# ---------------------------------------------------------------
$ws803 = $wb.Worksheets.Item("803")

$ws803.Cells.Item(36, 8).Value  = 3
$ws803.Cells.Item(54, 8).Value  = 2
$ws803.Cells.Item(99, 8).Value  = 3.5
$ws803.Cells.Item(117, 8).Value = 3.9
$ws803.Cells.Item(135, 8).Value = 3
$ws803.Cells.Item(153, 8).Value = 3
$ws803.Cells.Item(180, 8).Value = 2
$ws803.Cells.Item(207, 8).Value = 3.3
$ws803.Cells.Item(216, 8).Value = 2
$ws803.Cells.Item(261, 8).Value = 3
$ws803.Cells.Item(333, 8).Value = 1

# ---------------------------------------------------------------
# Sheet "115" — single-cell grade corrections (rows before the insert)
# ---------------------------------------------------------------
$ws115 = $wb.Worksheets.Item("115")

$ws115.Cells.Item(32, 8).Value  = 4
$ws115.Cells.Item(50, 8).Value  = 4.2
$ws115.Cells.Item(59, 8).Value  = 4.5
$ws115.Cells.Item(122, 8).Value = 4.5
$ws115.Cells.Item(131, 8).Value = 4
$ws115.Cells.Item(140, 8).Value = 4.3

# ---------------------------------------------------------------
# Sheet "115" — insert a new 9-row student block at row 164
# ---------------------------------------------------------------
# Use an existing, untouched 9-row block (rows 2-10) as a formatting /
# layout template: copying it and doing an "insert copied cells" at the
# destination both shifts the existing rows down by 9 AND seeds the new
# rows with the same A/E/F/G template values shared by every student.
$template = $ws115.Range("A2:H10")
$template.Copy()
$ws115.Range("A164:H172").Insert(-4121)  # xlShiftDown
$excel.CutCopyMode = $false

# Now overwrite the identity columns (B/C/D) and the grade column (H)
# for the freshly inserted rows with OSPINA ESCOBAR, SEBASTIAN's data.
# B/C are forced to text with a leading apostrophe so they keep storing
# as text (matching every other Matricula/DOCUMENTO cell in the sheet)
# instead of being auto-coerced to a number.
$newGrades = @(3, 3, 3.3, 1, 3, 3.8, 1, 1, 3)

for ($i = 0; $i -lt 9; $i++) {
    $r = 164 + $i
    $ws115.Cells.Item($r, 2).Value = "'240667"
    $ws115.Cells.Item($r, 3).Value = "'1035974995"
    $ws115.Cells.Item($r, 4).Value = "OSPINA ESCOBAR, SEBASTIAN"
    $ws115.Cells.Item($r, 8).Value = $newGrades[$i]
}
